$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15:E61").ClearContents()
$ws.Range("E6").Select()
